# This script updates the "cryptos" price table (rows 2-51) to reflect
# refreshed market data pulled on 2024-07-04 by the scheduled GitHub Actions job.
# Columns: A=index (unchanged), B=Coin, C=Link, D=Price, E=Volume(1h)
#
# Note: several "Price" values look like plain numbers (e.g. "532.39") but must
# stay stored as TEXT, exactly as they were originally (Excel would otherwise
# silently coerce them to numeric cells). We force text by switching the cell to
# the "@" (Text) number format before assigning the value, then restore the
# original "Normal" cell style so formatting/borders are unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Cells.Item(2,4).Value = "58.380.26"
$ws.Cells.Item(2,5).Value = "  -4.05%  "

# Row 3: Ethereum
$ws.Cells.Item(3,4).Value = "3.189.19"
$ws.Cells.Item(3,5).Value = "  -4.85%  "

# Row 4: TetherUSD
$ws.Cells.Item(4,5).Value = "  +0.04%  "

# Row 5: BNB
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "532.39"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = "  -5.93%  "

# Row 6: Solana
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "134.70"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = "  -7.84%  "

# Row 7: USDC
$ws.Cells.Item(7,5).Value = "  -0.06%  "

# Row 8: LidoStakedEther
$ws.Cells.Item(8,4).Value = "3.188.05"
$ws.Cells.Item(8,5).Value = "  -4.92%  "

# Row 9: XRP
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = "0.455"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value = "  -5.74%  "

# Row 10: Toncoin
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "7.35"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value = "  -6.87%  "

# Row 11: Dogecoin
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = "0.112"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = "  -7.59%  "

# Row 12: Cardano
$ws.Cells.Item(12,5).Value = "  -4.80%  "

# Row 13: WrappedliquidstakedEther2.0
$ws.Cells.Item(13,4).Value = "3.734.03"
$ws.Cells.Item(13,5).Value = "  -4.95%  "

# Row 14: TRON
$ws.Cells.Item(14,5).Value = "  -0.61%  "

# Row 15: Avalanche
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = "25.63"
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).Value = "  -8.01%  "

# Row 16: WrappedEther
$ws.Cells.Item(16,4).Value = "3.190.67"
$ws.Cells.Item(16,5).Value = "  -4.80%  "

# Row 17: WrappedBTC
$ws.Cells.Item(17,4).Value = "58.513.35"
$ws.Cells.Item(17,5).Value = "  -3.93%  "

# Row 18: ShibaInu
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = "0.0000155"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).Value = "  -8.09%  "

# Row 19: Polkadot
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = "5.87"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Value = "  -7.04%  "

# Row 20: Chainlink
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = "13.19"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value = "  -8.50%  "

# Row 21: Uniswap
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = "8.11"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Value = "  -8.60%  "

# Row 22: BitcoinCash
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = "358.63"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value = "  -4.65%  "

# Row 23: Dai
$ws.Cells.Item(23,5).Value = "  +0.02%  "

# Row 24: Litecoin
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "69.59"
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).Value = "  -7.07%  "

# Row 25: Polygon
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = "0.516"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value = "  -7.78%  "

# Row 26: WrappedeETH
$ws.Cells.Item(26,4).Value = "3.324.70"

# Row 27: Kaspa
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = "0.169"
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value = "  -2.82%  "

# Row 28: PEPE
$ws.Cells.Item(28,4).Value = "0.0₃0948"
$ws.Cells.Item(28,5).Value = "  -12.27%  "

# Row 29: Binance-PegBSC-USD
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = "0.997"
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).Value = "  -0.36%  "

# Row 30: RenderToken
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = "6.97"
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).Value = "  -5.50%  "

# Row 31: USDe
$ws.Cells.Item(31,5).Value = "  -0.10%  "

# Row 32: PancakeSwap
$ws.Cells.Item(32,5).Value = "  -8.24%  "

# Row 33: InternetComputer(DFINITY)
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = "7.00"
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).Value = "  -8.98%  "

# Row 34: EthereumClassic
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = "21.65"
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).Value = "  -5.18%  "

# Row 35: Fetch.AI
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = "1.20"
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).Value = "  -8.07%  "

# Row 36: NEARProtocol
$ws.Cells.Item(36,5).Value = "  -7.13%  "

# Row 37: Monero
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = "160.75"
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).Value = "  -4.76%  "

# Row 38: ImmutableX
$ws.Cells.Item(38,2).Value = "ImmutableX"
$ws.Cells.Item(38,3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = "1.43"
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).Value = "  -7.63%  "

# Row 39: Aptos
$ws.Cells.Item(39,2).Value = "Aptos"
$ws.Cells.Item(39,3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = "6.30"
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).Value = "  -7.61%  "

# Row 40: EnergySwap
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = "25.81"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value = "  -10.59%  "

# Row 41: Hedera
$ws.Cells.Item(41,5).Value = "  -6.09%  "

# Row 42: RenzoRestakedETH
$ws.Cells.Item(42,4).Value = "3.215.76"
$ws.Cells.Item(42,5).Value = "  -5.11%  "

# Row 43: OKB
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = "40.69"
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).Value = "  -3.77%  "

# Row 44: Mantle
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = "0.707"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).Value = "  -6.86%  "

# Row 45: Filecoin
$ws.Cells.Item(45,2).Value = "Filecoin"
$ws.Cells.Item(45,3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = "4.00"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value = "  -6.73%  "

# Row 46: ONDO
$ws.Cells.Item(46,2).Value = "ONDO"
$ws.Cells.Item(46,3).Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = "1.09"
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Value = "  -4.27%  "

# Row 47: Stacks
$ws.Cells.Item(47,5).Value = "  -7.64%  "

# Row 48: FirstDigitalUSD
$ws.Cells.Item(48,5).Value = "  -0.09%  "

# Row 49: Maker
$ws.Cells.Item(49,4).Value = "2.282.04"
$ws.Cells.Item(49,5).Value = "  -8.04%  "

# Row 50: Cosmos
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = "6.23"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value = "  -6.40%  "

# Row 51: InjectiveProtocol
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = "20.47"
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Value = "  -8.80%  "
